$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the profile id and human-readable name in row 2
# (Profile column A, Name column B) to the new "eICR Anonymized" naming.
$ws.Range("A2").Value = "eicr-anon-caculated-age"
$ws.Range("B2").Value = "eICR Anonymized Calculated Age"
